$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) contains values that look numeric (e.g. "334.14") but must
# remain plain text, matching the original inlineStr cell content. Temporarily
# force the column to Text format while writing values, then restore the style.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '28.534.32'
$ws.Range('E2').Value = '  -3.39%  '
$ws.Range('D3').Value = '1.850.09'
$ws.Range('E3').Value = '  -3.65%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.88%  '
$ws.Range('D5').Value = '334.14'
$ws.Range('E5').Value = '  +2.46%  '
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').Value = '0.4667'
$ws.Range('E7').Value = '  -3.14%  '
$ws.Range('D8').Value = '0.3924'
$ws.Range('E8').Value = '  -3.42%  '
$ws.Range('D9').Value = '46.44'
$ws.Range('E9').Value = '  -2.53%  '
$ws.Range('D10').Value = '0.07919'
$ws.Range('E10').Value = '  -3.91%  '
$ws.Range('D11').Value = '0.9854'
$ws.Range('E11').Value = '  -2.51%  '
$ws.Range('D12').Value = '22.24'
$ws.Range('E12').Value = '  -5.20%  '
$ws.Range('D13').Value = '2.019.73'
$ws.Range('E13').Value = '  +5.08%  '
$ws.Range('D14').Value = '5.852'
$ws.Range('E14').Value = '  -3.55%  '
$ws.Range('D15').Value = '7.018'
$ws.Range('E15').Value = '  -3.26%  '
$ws.Range('D16').Value = '0.06862'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '87.79'
$ws.Range('E17').Value = '  -4.22%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value = '1.003'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '0.00001008'
$ws.Range('E19').Value = '  -3.03%  '
$ws.Range('D20').Value = '17.07'
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').Value = '28.562.58'
$ws.Range('E22').Value = '  -3.35%  '
$ws.Range('D23').Value = '5.394'
$ws.Range('E23').Value = '  -5.12%  '
$ws.Range('D24').Value = '11.27'
$ws.Range('E24').Value = '  -5.14%  '
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.191.35'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '2.126'
$ws.Range('E26').Value = '  -2.68%  '
$ws.Range('D27').Value = '153.52'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').Value = '19.41'
$ws.Range('E28').Value = '  -3.09%  '
$ws.Range('D29').Value = '6.135'
$ws.Range('E29').Value = '  -6.03%  '
$ws.Range('D30').Value = '2.019'
$ws.Range('E30').Value = '  -3.79%  '
$ws.Range('D31').Value = '117.49'
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('D32').Value = '0.9786'
$ws.Range('E32').Value = '  -4.17%  '
$ws.Range('D33').Value = '0.09435'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('D34').Value = '5.369'
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('D35').Value = '3.498'
$ws.Range('E35').Value = '  -1.60%  '
$ws.Range('D36').Value = '1.347'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').Value = '0.06127'
$ws.Range('E37').Value = '  -3.67%  '
$ws.Range('D38').Value = '0.02200'
$ws.Range('E38').Value = '  -4.26%  '
$ws.Range('D39').Value = '1.162'
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('D40').Value = '0.5708'
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('D41').Value = '7.608'
$ws.Range('E41').Value = '  -3.99%  '
$ws.Range('D42').Value = '10.10'
$ws.Range('E42').Value = '  -6.04%  '
$ws.Range('D43').Value = '0.1794'
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('D44').Value = '2.399'
$ws.Range('E44').Value = '  -3.10%  '
$ws.Range('D45').Value = '1.223'
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('D46').Value = '11.88'
$ws.Range('E46').Value = '  -4.54%  '
$ws.Range('D47').Value = '0.5383'
$ws.Range('E47').Value = '  -3.37%  '
$ws.Range('D48').Value = '0.07158'
$ws.Range('E48').Value = '  -4.47%  '
$ws.Range('D49').Value = '1.904'
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('D50').Value = '113.84'
$ws.Range('E50').Value = '  -4.12%  '
$ws.Range('D51').Value = '42.95'
$ws.Range('E51').Value = '  +1.73%  '

$priceRange.Style = "Normal"
